$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Pneumonia / Cold columns)
$ws.Range("F2").Value = 0
$ws.Range("H2").Value = 1

# Row 5 (Pneumonia / Flu / Cold columns)
$ws.Range("F5").Value = 0.9
$ws.Range("G5").Value = 0.09
$ws.Range("H5").Value = 0

# Row 6 (Pneumonia / Cold columns)
$ws.Range("F6").Value = 0
$ws.Range("H6").Value = 1
